$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price/Volume columns remain text so numeric-looking strings
# (e.g. "0.9934", "44.93", dotted thousand-values) are not coerced to numbers/dates.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '28.136.76'
$ws.Range('E2').Value = '  -1.69%  '
$ws.Range('D3').Value = '1.814.02'
$ws.Range('E3').Value = '  +0.07%  '
$ws.Range('D4').Value = '0.9934'
$ws.Range('E4').Value = '  -0.91%  '
$ws.Range('D5').Value = '325.77'
$ws.Range('E5').Value = '  -1.02%  '
$ws.Range('D6').Value = '0.9898'
$ws.Range('E6').Value = '  -0.89%  '
$ws.Range('D7').Value = '0.4450'
$ws.Range('E7').Value = '  +0.28%  '
$ws.Range('D8').Value = '0.3785'
$ws.Range('E8').Value = '  -0.68%  '
$ws.Range('D9').Value = '44.93'
$ws.Range('E9').Value = '  +0.25%  '
$ws.Range('D10').Value = '0.07747'
$ws.Range('E10').Value = '  +1.69%  '
$ws.Range('D11').Value = '1.137'
$ws.Range('E11').Value = '  -1.53%  '
$ws.Range('B12').Value = 'Solana'
$ws.Range('C12').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D12').Value = '22.19'
$ws.Range('E12').Value = '  -3.13%  '
$ws.Range('B13').Value = 'BinanceUSD'
$ws.Range('C13').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D13').Value = '0.9902'
$ws.Range('E13').Value = '  -0.83%  '
$ws.Range('D14').Value = '6.264'
$ws.Range('E14').Value = '  -1.50%  '
$ws.Range('D15').Value = '7.506'
$ws.Range('E15').Value = '  -1.27%  '
$ws.Range('D16').Value = '1.794.48'
$ws.Range('E16').Value = '  -0.45%  '
$ws.Range('D17').Value = '91.89'
$ws.Range('E17').Value = '  +12.87%  '
$ws.Range('D18').Value = '0.00001081'
$ws.Range('E18').Value = '  -1.37%  '
$ws.Range('D19').Value = '0.06335'
$ws.Range('E19').Value = '  -6.13%  '
$ws.Range('D20').Value = '0.9912'
$ws.Range('E20').Value = '  -0.72%  '
$ws.Range('D21').Value = '17.57'
$ws.Range('E21').Value = '  -1.57%  '
$ws.Range('D22').Value = '6.299'
$ws.Range('E22').Value = '  -0.92%  '
$ws.Range('D23').Value = '0.5277'
$ws.Range('E23').Value = '  -2.90%  '
$ws.Range('D24').Value = '28.130.17'
$ws.Range('E24').Value = '  -1.71%  '
$ws.Range('D25').Value = '11.66'
$ws.Range('E25').Value = '  -2.34%  '
$ws.Range('D26').Value = '2.053'
$ws.Range('E26').Value = '  -15.05%  '
$ws.Range('D27').Value = '20.93'
$ws.Range('E27').Value = '  +0.73%  '
$ws.Range('D28').Value = '153.16'
$ws.Range('E28').Value = '  -0.18%  '
$ws.Range('D29').Value = '2.371'
$ws.Range('E29').Value = '  -0.66%  '
$ws.Range('D30').Value = '1.993.04'
$ws.Range('E30').Value = '  -0.70%  '
$ws.Range('D31').Value = '129.03'
$ws.Range('E31').Value = '  -3.37%  '
$ws.Range('D32').Value = '1.215'
$ws.Range('E32').Value = '  -4.76%  '
$ws.Range('D33').Value = '5.824'
$ws.Range('E33').Value = '  -1.12%  '
$ws.Range('D34').Value = '0.09232'
$ws.Range('E34').Value = '  -1.32%  '
$ws.Range('D35').Value = '3.645'
$ws.Range('E35').Value = '  -7.93%  '
$ws.Range('D36').Value = '12.82'
$ws.Range('E36').Value = '  +3.72%  '
$ws.Range('D37').Value = '0.02347'
$ws.Range('E37').Value = '  -0.22%  '
$ws.Range('D38').Value = '0.2191'
$ws.Range('E38').Value = '  -4.52%  '
$ws.Range('B39').Value = 'InternetComputer(DFINITY)'
$ws.Range('C39').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D39').Value = '5.165'
$ws.Range('E39').Value = '  -1.54%  '
$ws.Range('D40').Value = '0.6582'
$ws.Range('E40').Value = '  -1.37%  '
$ws.Range('B41').Value = 'Hedera'
$ws.Range('C41').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D41').Value = '0.06203'
$ws.Range('E41').Value = '  -3.11%  '
$ws.Range('D42').Value = '1.187'
$ws.Range('E42').Value = '  -2.11%  '
$ws.Range('D43').Value = '8.091'
$ws.Range('E43').Value = '  -1.34%  '
$ws.Range('D44').Value = '0.9900'
$ws.Range('E44').Value = '  -0.77%  '
$ws.Range('D45').Value = '1.389'
$ws.Range('E45').Value = '  -4.70%  '
$ws.Range('D46').Value = '13.83'
$ws.Range('E46').Value = '  -1.33%  '
$ws.Range('D47').Value = '0.6090'
$ws.Range('E47').Value = '  -0.98%  '
$ws.Range('D48').Value = '3.733'
$ws.Range('E48').Value = '  -2.34%  '
$ws.Range('D49').Value = '127.04'
$ws.Range('E49').Value = '  -1.96%  '
$ws.Range('D50').Value = '2.023'
$ws.Range('E50').Value = '  -2.09%  '
$ws.Range('B51').Value = 'Aave'
$ws.Range('C51').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D51').Value = '78.90'
$ws.Range('E51').Value = '  +0.29%  '
